# update report from 6-8 to 6-12
# Append new daily-report rows (2015-06-08 .. 2015-06-12) below the existing
# data, reusing the date-formatted style from the last existing row so no
# new number-format style gets created.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 19; Date = 42163; Text = "Learn the structure of Sharecare project." },
    @{ Row = 20; Date = 42164; Text = "Learn Solr and  the structure of Sharecare project." },
    @{ Row = 21; Date = 42165; Text = "Fix bug: Published question without an answer gives an OOPS page." },
    @{ Row = 22; Date = 42166; Text = "Fix bug: Published question without an answer gives an OOPS page." },
    @{ Row = 23; Date = 42167; Text = "Write unit test after fix the bug." }
)

foreach ($item in $newRows) {
    $r = $item.Row
    # Copy the format of the last pre-existing date cell (A18) so the new
    # cell picks up the same date number-format style instead of creating a
    # brand-new style entry.
    $ws.Range("A18").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false

    $ws.Range("A$r").Value2 = $item.Date
    $ws.Range("B$r").Value2 = $item.Text
}

# Match the author's final selection / active cell state.
$null = $ws.Range("B24").Select()
